# Update the build/version string throughout the workbook.
#
# Old version string: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# New version string: "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

# A2: "Version: <version string>"
$wsAbout.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended Citation containing the version string embedded in quotes.
$oldCitation = $wsAbout.Range("A6").Value()
$newCitation = $oldCitation.Replace($oldVersion, $newVersion)
$wsAbout.Range("A6").Value = $newCitation

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S (build_version) for data rows 2 through 10.
for ($row = 2; $row -le 10; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = 19
    $current = $cell.Value()
    if ($current -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
